$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap country names at rows 134/135 (Siria now before Guinea Ecuatorial in the
# underlying shared-string table; net effect: row 134 -> Siria, row 135 -> Guinea Ecuatorial)
$ws.Range("A134").Value = "Siria"
$ws.Range("A135").Value = "Guinea Ecuatorial"

# Update "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 23:00"

# Update numeric stats for affected countries
$ws.Range("B4").Value = 8380773
$ws.Range("C4").Value = 38108
$ws.Range("D4").Value = 5451197
$ws.Range("E4").Value = 2704933
$ws.Range("G4").Value = 361
$ws.Range("H4").Value = 224643
$ws.Range("B5").Value = 7547759
$ws.Range("C5").Value = 55032
$ws.Range("D5").Value = 6658937
$ws.Range("E5").Value = 774193
$ws.Range("G5").Value = 565
$ws.Range("H5").Value = 114629
$ws.Range("B21").Value = 366944
$ws.Range("C21").Value = 5211
$ws.Range("E21").Value = 67078
$ws.Range("B117").Value = 7752
$ws.Range("C117").Value = 114
$ws.Range("E117").Value = 1193
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 86
$ws.Range("B129").Value = 5297
$ws.Range("C129").Value = 16
$ws.Range("D129").Value = 3652
$ws.Range("E129").Value = 1549
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 96
$ws.Range("B134").Value = 5077
$ws.Range("C134").Value = 44
$ws.Range("D134").Value = 1528
$ws.Range("E134").Value = 3301
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 248
$ws.Range("B135").Value = 5070
$ws.Range("D135").Value = 4954
$ws.Range("E135").Value = 33
$ws.Range("H135").Value = 83
$ws.Range("B136").Value = 4974
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 4783
$ws.Range("E136").Value = 157
$ws.Range("B146").Value = 3734
$ws.Range("C146").Value = 24
$ws.Range("D146").Value = 2654
$ws.Range("E146").Value = 971
$ws.Range("B161").Value = 2056
$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 1338
$ws.Range("E161").Value = 121
$ws.Range("G161").Value = 1
$ws.Range("H161").Value = 597
$ws.Range("B169").Value = 933
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 898
$ws.Range("E169").Value = 20
